# "saving all the work"
#
# This script:
#  - Cleans up Sheet1 ("Sheet"): drops the extra duplicate rows (3-8),
#    removes the mailto hyperlinks + their underline "Hyperlink" style,
#    fixes up the selection, and rewrites the header/data to the new
#    values (id / name / email, 99003758 / sravan / abc).
#  - Renames "Master_sheet" to "Master6" and adds the new blank
#    Master5..Master1/Master sheets after it, matching the sheet order
#    and ids in the workbook.

$wb = $excel.ActiveWorkbook

# --- Sheet1 ("Sheet") edits ---------------------------------------------
$ws = $wb.Worksheets.Item("Sheet")

# Drop the duplicate rows, leaving only the header row and one data row.
$ws.Rows("3:8").Delete()

# Remove the mailto: hyperlinks that were on column C.
$ws.Hyperlinks.Delete()

# Put C2 back on the plain "Normal" style (removes the underlined
# hyperlink-colored font that used to be applied to it).
$ws.Range("C2").Style = "Normal"

# New header row text.
$ws.Range("A1").Value = "id "
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "email"

# New data row values.
$ws.Range("A2").Value = 99003758
$ws.Range("B2").Value = "sravan"
$ws.Range("C2").Value = "abc"

# The named "Hyperlink" cell style is no longer used anywhere - drop it.
$wb.Styles.Item("Hyperlink").Delete()

# Selection should just be the single cell C2 now.
$ws.Range("C2").Select()

# --- Workbook sheet restructuring ---------------------------------------
$master = $wb.Worksheets.Item("Master_sheet")
$master.Name = "Master6"

$newNames = @("Master5", "Master4", "Master3", "Master2", "Master1", "Master")
$prev = $master
foreach ($n in $newNames) {
    $newSheet = $wb.Worksheets.Add($null, $prev)
    $newSheet.Name = $n
    $prev = $newSheet
}

# Keep "Sheet" as the active/selected tab.
$ws.Activate()
